# Insert a new column before column F on every sheet. Excel's native
# "insert column" behaviour:
#   - shifts the old column F (and everything to its right) one column to
#     the right, so the old F..O data lands in G..P (formulas get their
#     column references auto-adjusted, e.g. O30 -> P30), and
#   - the newly inserted column F inherits its formatting from the column
#     immediately to its left (column E), i.e. the plain "label" style,
#     leaving it blank.
# That is exactly the edit described in the commit: "Make sure UOM's
# right is blank" (the new blank UOM column) while preserving the rest of
# the numeric columns (now shifted out to include a new rightmost column).

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Inputs", "Balance Sheet", "Corkscrew")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F1").EntireColumn.Insert()
}
